$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 520
$ws1.Range("F8").Value = 851
$ws1.Range("F9").Value = 568
$ws1.Range("F11").Value = 346
$ws1.Range("F14").Value = 1096
$ws1.Range("F17").Value = 440
$ws1.Range("F18").Value = 124
$ws1.Range("F19").Value = 268
$ws1.Range("F22").Value = 511
$ws1.Range("F25").Value = 407

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 175

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 520
$ws4.Range("F13").Value = 851
$ws4.Range("F14").Value = 568
$ws4.Range("F16").Value = 346
$ws4.Range("F19").Value = 1096
$ws4.Range("F24").Value = 440
$ws4.Range("F26").Value = 124
$ws4.Range("F28").Value = 268
$ws4.Range("F31").Value = 175
$ws4.Range("F33").Value = 511
$ws4.Range("F38").Value = 407
